$d = $word.ActiveDocument

$replacements = @(
    @("230÷9=", "296÷5="),
    @("316÷2=", "868÷2="),
    @("266÷9=", "501÷6="),
    @("856÷7=", "361÷7="),
    @("233÷8=", "159÷9="),
    @("910÷3=", "585÷9="),
    @("558÷9=", "338÷2="),
    @("421÷8=", "815÷7="),
    @("194÷3=", "590÷7="),
    @("782÷6=", "393÷7="),
    @("222÷2=", "992÷5="),
    @("582÷6=", "750÷6="),
    @("522÷3=", "890÷8="),
    @("108÷9=", "106÷7="),
    @("280÷9=", "870÷4="),
    @("759÷8=", "940÷6="),
    @("334÷5=", "866÷7="),
    @("871÷6=", "653÷6="),
    @("133÷7=", "477÷9="),
    @("964÷4=", "883÷3="),
    @("123÷9=", "735÷6="),
    @("483÷3=", "552÷9="),
    @("391÷8=", "914÷7="),
    @("346÷4=", "577÷4="),
    @("629÷7=", "692÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying $($replacements.Count) replacements"
